$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting of the last existing data row (781) down to each new row
# (Range.Copy-with-destination preserves the exact cell style indices used by
#  row 781, so we don't create redundant style entries in styles.xml).
$ws.Range("A781:D781").Copy($ws.Range("A782:D782"))
$ws.Range("A781:D781").Copy($ws.Range("A783:D783"))
$ws.Range("A781:D781").Copy($ws.Range("A784:D784"))
$ws.Range("A781:D781").Copy($ws.Range("A785:D785"))
$ws.Range("A781:D781").Copy($ws.Range("A786:D786"))
$ws.Range("A781:D781").Copy($ws.Range("A787:D787"))
$ws.Range("A781:D781").Copy($ws.Range("A788:D788"))
$ws.Range("A781:D781").Copy($ws.Range("A789:D789"))
$ws.Range("A781:D781").Copy($ws.Range("A790:D790"))
$ws.Range("A781:D781").Copy($ws.Range("A791:D791"))
$ws.Range("A781:D781").Copy($ws.Range("A792:D792"))
$ws.Range("A781:D781").Copy($ws.Range("A793:D793"))
$ws.Range("A781:D781").Copy($ws.Range("A794:D794"))
$ws.Range("A781:D781").Copy($ws.Range("A795:D795"))
$ws.Range("A781:D781").Copy($ws.Range("A796:D796"))
$ws.Range("A781:D781").Copy($ws.Range("A797:D797"))
$ws.Range("A781:D781").Copy($ws.Range("A798:D798"))
$ws.Range("A781:D781").Copy($ws.Range("A799:D799"))
$ws.Range("A781:D781").Copy($ws.Range("A800:D800"))

# Now overwrite with the real values for each new row.
$ws.Cells.Item(782, 1).Value = "2022-05-27"
$ws.Cells.Item(782, 2).Value = 70
$ws.Cells.Item(782, 3).Value = 70
$ws.Cells.Item(782, 4).Formula = "=C782/B782"
$ws.Cells.Item(783, 1).Value = "2022-05-28"
$ws.Cells.Item(783, 2).Value = 50
$ws.Cells.Item(783, 3).Value = 48
$ws.Cells.Item(783, 4).Formula = "=C783/B783"
$ws.Cells.Item(784, 1).Value = "2022-05-29"
$ws.Cells.Item(784, 2).Value = 49
$ws.Cells.Item(784, 3).Value = 47
$ws.Cells.Item(784, 4).Formula = "=C784/B784"
$ws.Cells.Item(785, 1).Value = "2022-05-30"
$ws.Cells.Item(785, 2).Value = 56
$ws.Cells.Item(785, 3).Value = 56
$ws.Cells.Item(785, 4).Formula = "=C785/B785"
$ws.Cells.Item(786, 1).Value = "2022-05-31"
$ws.Cells.Item(786, 2).Value = 51
$ws.Cells.Item(786, 3).Value = 50
$ws.Cells.Item(786, 4).Formula = "=C786/B786"
$ws.Cells.Item(787, 1).Value = "2022-06-01"
$ws.Cells.Item(787, 2).Value = 64
$ws.Cells.Item(787, 3).Value = 58
$ws.Cells.Item(787, 4).Formula = "=C787/B787"
$ws.Cells.Item(788, 1).Value = "2022-06-02"
$ws.Cells.Item(788, 2).Value = 73
$ws.Cells.Item(788, 3).Value = 67
$ws.Cells.Item(788, 4).Formula = "=C788/B788"
$ws.Cells.Item(789, 1).Value = "2022-06-03"
$ws.Cells.Item(789, 2).Value = 84
$ws.Cells.Item(789, 3).Value = 73
$ws.Cells.Item(789, 4).Formula = "=C789/B789"
$ws.Cells.Item(790, 1).Value = "2022-06-04"
$ws.Cells.Item(790, 2).Value = 59
$ws.Cells.Item(790, 3).Value = 57
$ws.Cells.Item(790, 4).Formula = "=C790/B790"
$ws.Cells.Item(791, 1).Value = "2022-06-05"
$ws.Cells.Item(791, 2).Value = 64
$ws.Cells.Item(791, 3).Value = 59
$ws.Cells.Item(791, 4).Formula = "=C791/B791"
$ws.Cells.Item(792, 1).Value = "2022-06-06"
$ws.Cells.Item(792, 2).Value = 70
$ws.Cells.Item(792, 3).Value = 69
$ws.Cells.Item(792, 4).Formula = "=C792/B792"
$ws.Cells.Item(793, 1).Value = "2022-06-07"
$ws.Cells.Item(793, 2).Value = 66
$ws.Cells.Item(793, 3).Value = 62
$ws.Cells.Item(793, 4).Formula = "=C793/B793"
$ws.Cells.Item(794, 1).Value = "2022-06-08"
$ws.Cells.Item(794, 2).Value = 78
$ws.Cells.Item(794, 3).Value = 78
$ws.Cells.Item(794, 4).Formula = "=C794/B794"
$ws.Cells.Item(795, 1).Value = "2022-06-09"
$ws.Cells.Item(795, 2).Value = 80
$ws.Cells.Item(795, 3).Value = 76
$ws.Cells.Item(795, 4).Formula = "=C795/B795"
$ws.Cells.Item(796, 1).Value = "2022-06-10"
$ws.Cells.Item(796, 2).Value = 70
$ws.Cells.Item(796, 3).Value = 67
$ws.Cells.Item(796, 4).Formula = "=C796/B796"
$ws.Cells.Item(797, 1).Value = "2022-06-11"
$ws.Cells.Item(797, 2).Value = 60
$ws.Cells.Item(797, 3).Value = 58
$ws.Cells.Item(797, 4).Formula = "=C797/B797"
$ws.Cells.Item(798, 1).Value = "2022-06-12"
$ws.Cells.Item(798, 2).Value = 63
$ws.Cells.Item(798, 3).Value = 61
$ws.Cells.Item(798, 4).Formula = "=C798/B798"
$ws.Cells.Item(799, 1).Value = "2022-06-13"
$ws.Cells.Item(799, 2).Value = 86
$ws.Cells.Item(799, 3).Value = 81
$ws.Cells.Item(799, 4).Formula = "=C799/B799"
$ws.Cells.Item(800, 1).Value = "2022-06-14"
$ws.Cells.Item(800, 2).Value = 79
$ws.Cells.Item(800, 3).Value = 74
$ws.Cells.Item(800, 4).Formula = "=C800/B800"

# Move the visible selection to match the saved view state.
$ws.Range("F799").Select() | Out-Null

